$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bold the existing "Company" column cells (rows 2-18) -> introduces the
# new cellXfs entry (fontId=1, applyFont=1, no fill) used as style index 2.
$ws.Range("A2:A18").Font.Bold = $true

# New data rows 19-24
$newRows = @(
    @("DDD", "Aasdfsdad", "fsdfssdsdf", "sdfsadsdaf", 12),
    @("EEE", "Aasdfsdad", "fsdfssdsdf", "sdfsadsdaf", 12),
    @("ZZZZZZZZZZ", "Aasdfsdad", "fsdfssdsdf", "sdfsadsdaf", 12),
    @("AAAAAAAAAAAAAA", "XXXXXXXXXXXXXX", "XXXXXXXXXXXXXX", "XXXXXXXXXXXX", 111),
    @("XCVSAMDVSDJKLANFKLJASDNFLJKSNADL", "XXXXXXXXXXXXXX", "XXXXXXXXXXXXXX", "XXXXXXXXXXXX", 111),
    @("XCVSAMDVSDJKLANFK", "XXXXXXXXXXXXXX", "XXXXXXXXXXXXXX", "XXXXXXXXXXXX", 111)
)

$r = 19
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
    $ws.Range("A$r").Font.Bold = $true
    $r = $r + 1
}
